# Add the answer to sub-question (a) "What are the constraints?" in the
# second "Break the problem apart" section (the sock-drawer problem).
#
# The target is the paragraph that immediately follows the second
# occurrence of "a) What are the constraints?" -- it is an otherwise
# empty paragraph that only contains the _GoBack bookmark.

$d = $word.ActiveDocument

$targetIndex = -1
$occurrence = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "a\) What are the constraints\?") {
        $occurrence = $occurrence + 1
        if ($occurrence -eq 2) {
            $targetIndex = $i + 1
            break
        }
    }
}

$target = $d.Paragraphs.Item($targetIndex)

# First run: insert right at the start of the (empty) paragraph, ahead
# of the _GoBack bookmark that already lives there.
$r = $target.Range
$r.Collapse(1)
$r.InsertBefore("The constraints for this problem would definitely doing it in the dark because it not being in the dark you could just collect the socks you desire. Another constraint")
$r.Font.Size = 11

# Second run: appended right after the first one. Toggling a property
# and back keeps this text in its own run instead of merging back into
# the previous run (matching the two <w:r> elements in the target).
$r.Collapse(0)
$r.InsertAfter(" would definitely be the odd number of socks indicating that there is going to be an extra black and brown sock.")
$r.Font.Bold = $true
$r.Font.Bold = $false
$r.Font.Size = 11
